$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-5) and make room for the new 2-10 data block
$ws.Range("A2:T10").Clear()

# Column A
$ws.Range("A2").Value = "ECs"
$ws.Range("A3").Value = "ECs"
$ws.Range("A4").Value = "ECs"
$ws.Range("A5").Value = "FAPs"
$ws.Range("A6").Value = "FAPs"
$ws.Range("A7").Value = "FAPs"
$ws.Range("A8").Value = "sCs"
$ws.Range("A9").Value = "sCs"
$ws.Range("A10").Value = "sCs"

# Column B
$ws.Range("B2").Value = "Dcn"
$ws.Range("B3").Value = "Dcn"
$ws.Range("B4").Value = "Dcn"
$ws.Range("B5").Value = "Dcn"
$ws.Range("B6").Value = "Dcn"
$ws.Range("B7").Value = "Dcn"
$ws.Range("B8").Value = "Dcn"
$ws.Range("B9").Value = "Dcn"
$ws.Range("B10").Value = "Dcn"

# Column C
$ws.Range("C2").Value = "Met"
$ws.Range("C3").Value = "Met"
$ws.Range("C4").Value = "Met"
$ws.Range("C5").Value = "Met"
$ws.Range("C6").Value = "Met"
$ws.Range("C7").Value = "Met"
$ws.Range("C8").Value = "Met"
$ws.Range("C9").Value = "Met"
$ws.Range("C10").Value = "Met"

# Column D
$ws.Range("D2").Value = "ECs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("D4").Value = "sCs"
$ws.Range("D5").Value = "ECs"
$ws.Range("D6").Value = "FAPs"
$ws.Range("D7").Value = "sCs"
$ws.Range("D8").Value = "ECs"
$ws.Range("D9").Value = "FAPs"
$ws.Range("D10").Value = "sCs"

# Column E
$ws.Range("E2").Value = 2
$ws.Range("E3").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("E5").Value = 3
$ws.Range("E6").Value = 3
$ws.Range("E7").Value = 3
$ws.Range("E8").Value = 3
$ws.Range("E9").Value = 3
$ws.Range("E10").Value = 3

# Column F
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("F5").Value = 1
$ws.Range("F6").Value = 1
$ws.Range("F7").Value = 1
$ws.Range("F8").Value = 1
$ws.Range("F9").Value = 1
$ws.Range("F10").Value = 1

# Column G
$ws.Range("G2").Value = 6.615074666666668
$ws.Range("G3").Value = 6.615074666666668
$ws.Range("G4").Value = 6.615074666666668
$ws.Range("G5").Value = 7285.701009
$ws.Range("G6").Value = 7285.701009
$ws.Range("G7").Value = 7285.701009
$ws.Range("G8").Value = 126.6246363333333
$ws.Range("G9").Value = 126.6246363333333
$ws.Range("G10").Value = 126.6246363333333

# Column H
$ws.Range("H2").Value = 19.845224
$ws.Range("H3").Value = 19.845224
$ws.Range("H4").Value = 19.845224
$ws.Range("H5").Value = 21857.103027
$ws.Range("H6").Value = 21857.103027
$ws.Range("H7").Value = 21857.103027
$ws.Range("H8").Value = 379.873909
$ws.Range("H9").Value = 379.873909
$ws.Range("H10").Value = 379.873909

# Column I
$ws.Range("I2").Value = 0.0008916467884469992
$ws.Range("I3").Value = 0.0008916467884469992
$ws.Range("I4").Value = 0.0008916467884469992
$ws.Range("I5").Value = 0.9820406017477925
$ws.Range("I6").Value = 0.9820406017477925
$ws.Range("I7").Value = 0.9820406017477925
$ws.Range("I8").Value = 0.01706775146376063
$ws.Range("I9").Value = 0.01706775146376063
$ws.Range("I10").Value = 0.01706775146376063

# Column J
$ws.Range("J2").Value = 0.0008916467884469989
$ws.Range("J3").Value = 0.0008916467884469989
$ws.Range("J4").Value = 0.0008916467884469989
$ws.Range("J5").Value = 0.9820406017477923
$ws.Range("J6").Value = 0.9820406017477923
$ws.Range("J7").Value = 0.9820406017477923
$ws.Range("J8").Value = 0.01706775146376063
$ws.Range("J9").Value = 0.01706775146376063
$ws.Range("J10").Value = 0.01706775146376063

# Column K
$ws.Range("K2").Value = 3
$ws.Range("K3").Value = 1
$ws.Range("K4").Value = 3
$ws.Range("K5").Value = 3
$ws.Range("K6").Value = 1
$ws.Range("K7").Value = 3
$ws.Range("K8").Value = 3
$ws.Range("K9").Value = 1
$ws.Range("K10").Value = 3

# Column L
$ws.Range("L2").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("L4").Value = 1
$ws.Range("L5").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("L7").Value = 1
$ws.Range("L8").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("L10").Value = 1

# Column M
$ws.Range("M2").Value = 2.254050666666667
$ws.Range("M3").Value = 0.05371366666666667
$ws.Range("M4").Value = 1.605544
$ws.Range("M5").Value = 2.254050666666667
$ws.Range("M6").Value = 0.05371366666666667
$ws.Range("M7").Value = 1.605544
$ws.Range("M8").Value = 2.254050666666667
$ws.Range("M9").Value = 0.05371366666666667
$ws.Range("M10").Value = 1.605544

# Column N
$ws.Range("N2").Value = 6.762152
$ws.Range("N3").Value = 0.161141
$ws.Range("N4").Value = 4.816632
$ws.Range("N5").Value = 6.762152
$ws.Range("N6").Value = 0.161141
$ws.Range("N7").Value = 4.816632
$ws.Range("N8").Value = 6.762152
$ws.Range("N9").Value = 0.161141
$ws.Range("N10").Value = 4.816632

# Column O
$ws.Range("O2").Value = 0.5759961839619929
$ws.Range("O3").Value = 0.01372589688605336
$ws.Range("O4").Value = 0.4102779191519537
$ws.Range("O5").Value = 0.5759961839619929
$ws.Range("O6").Value = 0.01372589688605336
$ws.Range("O7").Value = 0.4102779191519537
$ws.Range("O8").Value = 0.5759961839619929
$ws.Range("O9").Value = 0.01372589688605336
$ws.Range("O10").Value = 0.4102779191519537

# Column P
$ws.Range("P2").Value = 0.575996183961993
$ws.Range("P3").Value = 0.01372589688605336
$ws.Range("P4").Value = 0.4102779191519537
$ws.Range("P5").Value = 0.575996183961993
$ws.Range("P6").Value = 0.01372589688605336
$ws.Range("P7").Value = 0.4102779191519537
$ws.Range("P8").Value = 0.575996183961993
$ws.Range("P9").Value = 0.01372589688605336
$ws.Range("P10").Value = 0.4102779191519537

# Column Q
$ws.Range("Q2").Value = 14.91071346244978
$ws.Range("Q3").Value = 0.3553199156204445
$ws.Range("Q4").Value = 10.62079344061867
$ws.Range("Q5").Value = 16422.33921647046
$ws.Range("Q6").Value = 391.341715430423
$ws.Range("Q7").Value = 11697.5135407939
$ws.Range("Q8").Value = 285.4183459435743
$ws.Range("Q9").Value = 6.801473507796556
$ws.Range("Q10").Value = 203.3014251171654

# Column R
$ws.Range("R2").Value = 134.196421162048
$ws.Range("R3").Value = 3.197879240584
$ws.Range("R4").Value = 95.58714096556801
$ws.Range("R5").Value = 147801.0529482341
$ws.Range("R6").Value = 3522.075438873807
$ws.Range("R7").Value = 105277.6218671451
$ws.Range("R8").Value = 2568.765113492168
$ws.Range("R9").Value = 61.21326157016901
$ws.Range("R10").Value = 1829.712826054488

# Column S
$ws.Range("S2").Value = 0.0005135851475874379
$ws.Range("S3").Value = 0.00001223865187700415
$ws.Range("S4").Value = 0.0003658229889825571
$ws.Range("S5").Value = 0.5656516391024677
$ws.Range("S6").Value = 0.01347938803750799
$ws.Range("S7").Value = 0.4029095746078167
$ws.Range("S8").Value = 0.009830959711937842
$ws.Range("S9").Value = 0.0002342701966683647
$ws.Range("S10").Value = 0.007002521555154423

# Column T
$ws.Range("T2").Value = 0.0005135851475874379
$ws.Range("T3").Value = 0.00001223865187700414
$ws.Range("T4").Value = 0.000365822988982557
$ws.Range("T5").Value = 0.5656516391024677
$ws.Range("T6").Value = 0.01347938803750799
$ws.Range("T7").Value = 0.4029095746078167
$ws.Range("T8").Value = 0.009830959711937842
$ws.Range("T9").Value = 0.0002342701966683647
$ws.Range("T10").Value = 0.007002521555154422

